# Apply "change to the databases" edits to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set new cell values in the same order the shared strings were introduced
# so the underlying sharedStrings table layout matches the authored edit.
$ws.Range("A21").Value = "preference"
$ws.Range("B1").Value = "id"
$ws.Range("A9").Value = "username id "
$ws.Range("E7").Value = "exams table"
$ws.Range("E8").Value = "exam id"
$ws.Range("E9").Value = "exam name"
$ws.Range("E10").Value = "qualification id"
$ws.Range("E11").Value = "qualification name"
$ws.Range("A8").Value = "request tables"
$ws.Range("E16").Value = "response table"
$ws.Range("B8").Value = "request id"
$ws.Range("E17").Value = "reqid"
$ws.Range("E18").Value = "vol id"
$ws.Range("E19").Value = "date"

# Re-use of existing strings
$ws.Range("B6").Value = "image"

# row 13 ("image") is no longer needed on its own - value now lives in B6
$ws.Range("A13").Value = ""

# Update sheet view: drop the frozen/top-left cell and move the active selection
$ws.Range("E20").Select()

Write-Output "edits applied"
